$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Inserting a column shifts the old N/O/P columns one to the right (-> O/P/Q),
# carrying over their values/styles, and leaves a new blank column N behind
# (Excel copies the formatting of the column to its left, i.e. column M).
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and select the cell to the
# right of the newly inserted column on the last data row.
$ws.Activate()
$ws.Range("R6").Select()
